$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = [double]"22.42000000000007"
$ws.Range("G2").Value = [double]"1.884714606603666e-12"
$ws.Range("H2").Value = [double]"5.413341371674005e-12"
$ws.Range("K2").Value = [double]"44.54479696061006"
$ws.Range("L2").Value = "[31.351176566699593, 57.738417354520536]"
$ws.Range("M2").Value = [double]"2.900701900898639e-10"
$ws.Range("N2").Value = [double]"2.900701900898639e-10"
$ws.Range("O2").Value = [double]"1.717026615475501"
$ws.Range("P2").Value = "[1.3899739268135, 2.0440793041375027]"
$ws.Range("S2").Value = [double]"63.50822187369947"
$ws.Range("T2").Value = "[55.601269892907496, 71.41517385449144]"
$ws.Range("W2").Value = [double]"16.29321321321326"
$ws.Range("X2").Value = [double]"15.12620620620625"
$ws.Range("Y2").Value = [double]"17.46022022022028"

# Row 3 updates
$ws.Range("E3").Value = [double]"24.72000000000043"
$ws.Range("H3").Value = [double]"3.566981605221386e-16"
$ws.Range("K3").Value = [double]"50.39240615471194"
$ws.Range("L3").Value = "[40.42894146600064, 60.35587084342323]"
$ws.Range("O3").Value = [double]"-2.566105711040311"
$ws.Range("P3").Value = "[-2.7673689040630816, -2.364842518017541]"
$ws.Range("S3").Value = [double]"65.25210718573793"
$ws.Range("T3").Value = "[60.02104237168578, 70.48317199979007]"
$ws.Range("W3").Value = [double]"10.09585585585603"
$ws.Range("X3").Value = [double]"9.304024024024182"
$ws.Range("Y3").Value = [double]"10.88768768768788"
